# Update Name of Algo
# Apply updated numeric results to Sheet1 (result_data_RandomForest)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A11").Value = -21.74250000000002
$ws.Range("D11").Value = -7.51329999999999
$ws.Range("A12").Value = -21.3892
$ws.Range("A15").Value = -21.96640000000001
$ws.Range("D23").Value = -8.305
$ws.Range("A27").Value = -21.97069999999999
$ws.Range("A28").Value = -21.9355
$ws.Range("D28").Value = -8.157999999999994
$ws.Range("A31").Value = -21.8804
$ws.Range("A32").Value = -21.58440000000001
$ws.Range("D32").Value = -7.110299999999993
$ws.Range("D34").Value = -8.005699999999997
$ws.Range("A36").Value = -20.5421
$ws.Range("D36").Value = -6.860900000000003
$ws.Range("D37").Value = -8.424700000000005
$ws.Range("A38").Value = -19.61859999999998
$ws.Range("D42").Value = -8.809999999999995
$ws.Range("A46").Value = -21.7743
$ws.Range("D49").Value = -7.956799999999999
$ws.Range("A54").Value = -21.94050000000001
$ws.Range("D54").Value = -7.915999999999997
$ws.Range("A55").Value = -22.18870000000001
$ws.Range("A56").Value = -22.14170000000001
$ws.Range("A67").Value = -21.46969999999997
$ws.Range("A69").Value = -21.65089999999997
$ws.Range("A72").Value = -21.8435
$ws.Range("A73").Value = -20.05169999999999
$ws.Range("D78").Value = -8.058600000000002
$ws.Range("D80").Value = -8.074100000000001
$ws.Range("A83").Value = -21.65969999999999
$ws.Range("A86").Value = -21.5981
$ws.Range("A91").Value = -20.53249999999998
$ws.Range("A93").Value = -21.35890000000001
$ws.Range("D97").Value = -8.033199999999994
$ws.Range("A99").Value = -21.85020000000002
$ws.Range("D99").Value = -7.884699999999994
$ws.Range("D100").Value = -8.17799999999999
$ws.Range("D101").Value = -8.280699999999994
$ws.Range("A104").Value = -21.6141
$ws.Range("A105").Value = -19.67069999999999
